# pin-allocation.xlsx edit:
# "allocated pins to the rest of the features."
#
# Sheet1 ("Pin Name" / i-o requirements sheet):
#   - the "step/dir/enable + protect" placeholder pin (F3) is no longer needed
#   - the two outstanding "assign pins" / "assign protect pins" status markers
#     (red text) are now complete -> flip them to the same "ok" (green) marker
#     used by every other finished row
#   - the "assign pins" line item in the "more todo:" list is done -> remove it
#
# Sheet2 (pin allocation table):
#   - fill in the "My Usage" column for the previously-unassigned protect /
#     buzzer / spindle-relay pins (column K, rows 7-15)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet1
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# F3 ("P0.2" / step-dir-enable+protect leftover note) is no longer used.
$ws1.Range("F3").Clear()

# E3: "assign protect pins" (red, left aligned) -> "ok" (green, left aligned)
$ws1.Range("E3").Style = "Normal"
$ws1.Range("E3").Value = "ok"
$ws1.Range("E3").Font.Color = 32768
$ws1.Range("E3").HorizontalAlignment = -4131

# E12: "assign pins" (red) -> "ok" (green)
$ws1.Range("E12").Style = "Normal"
$ws1.Range("E12").Value = "ok"
$ws1.Range("E12").Font.Color = 32768

# Remove the now-finished "assign pins" to-do row.
$ws1.Range("A26").Clear()

$ws1.Activate()
$ws1.Range("E12").Select()

# ---------------------------------------------------------------------------
# Sheet2
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("K7").Value = "x protect"
$ws2.Range("K8").Value = "y protect"
$ws2.Range("K9").Value = "z protect"
$ws2.Range("K11").Value = "a protect"
$ws2.Range("K12").Value = "b protect"
$ws2.Range("K14").Value = "BUZZER"
$ws2.Range("K15").Value = "SPINDLE RELAY"

$krange = $ws2.Range("K7,K8,K9,K11,K12,K14,K15")
$krange.Font.Color = 32768

$ws2.Columns.Item(11).ColumnWidth = 14.33203125

$ws2.Activate()
$ws2.Range("K15").Select()
